$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replicate the formatting (styles) of the last existing bibliography row
# onto the four new rows before filling in content. Copy only the populated
# (non-contiguous) column groups so blank columns E/H aren't materialised.
$ws.Range("A118:D118").Copy()
$ws.Range("A119:D122").PasteSpecial(-4122)
$ws.Range("F118:G118").Copy()
$ws.Range("F119:G122").PasteSpecial(-4122)
$ws.Range("I118:K118").Copy()
$ws.Range("I119:K122").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

for ($row = 119; $row -le 122; $row++) {
    $ws.Rows.Item($row).RowHeight = 17
}

# New bibliography rows (body mass & competition), appended after row 118.
$rows = @(
    @{ Row = 119; Title = "Assessment of energy reserves by damselflies engaged in aerial contests for mating territories"; Journal = "Anim. Behav."; Year = 1994; Author = "Marden & Rollins"; Point = "males *Calopteryx maculata* with more fat win contests for mating territory in 88 % of cases"; PointFirst = $false },
    @{ Row = 120; Title = "Female competition and its evolutionary consequences in mammals"; Journal = "Biological Reviews"; Year = 2011; Author = "Stockley & Bro-Jørgensen"; Point = '"there is evidence that female dominance may often be correlated with age or body size"; "larger females consistently dominate" (e.g., feral ponies and african elephants)'; PointFirst = $false },
    @{ Row = 121; Title = "Dominance, aggression frequencies and modes of aggressive competition in feral pony mares"; Journal = "Anim. Behav."; Year = 1990; Author = "Rutberg & Greenberg"; Point = "Older and Larger Equus caballus are dominant"; PointFirst = $false },
    @{ Row = 122; Title = "Dominance rank relationships among wild female African elephants, Loxodonta africana"; Journal = "Anim. Behav."; Year = 2006; Author = "Archie et al."; Point = "Older and Larger Loxodonta africana are dominant"; PointFirst = $true }
)

foreach ($r in $rows) {
    $row = $r.Row

    if ($r.PointFirst) {
        $ws.Range("G$row").Value = $r.Point
    }

    $ws.Range("A$row").Value = $r.Title
    $ws.Range("B$row").Value = $r.Journal
    $ws.Range("C$row").Value = $r.Year
    $ws.Range("D$row").Value = $r.Author
    $ws.Range("F$row").Value = "≈"

    if (-not $r.PointFirst) {
        $ws.Range("G$row").Value = $r.Point
    }

    $ws.Range("I$row").Value = "yes"
    $ws.Range("J$row").Value = "yes"
    $ws.Range("K$row").Value = "yes"
}

# Scroll / selection state reflecting the end of the edit session.
$ws.Application.ActiveWindow.ScrollRow = 100
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B122").Select()
